$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (header row 1 stays unchanged).
# Row order changed: Pelayanan, Kestabilan, Harga (Kecepatan row removed)
$data = @(
    @("Pelayanan", 9, 1, 0, 1, 1, 0.5000000000000001),
    @("Kestabilan", 8, 2, 0.5, 1.5, 0.6666666666666666, 0.3333333333333334),
    @("Harga", 6, 3, 1, 2, 0.3333333333333333, 0.1666666666666667)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# Remove the now-obsolete 5th row (previously "Harga", now no longer needed
# since data only spans rows 1-4).
$ws.Rows.Item(5).Delete()
